$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.155.68"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").Value = "2.174.74"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("E5").Value = "  -2.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.613"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.19"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -5.38%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -6.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.97"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -9.12%  "
$ws.Range("E11").Value = "  -3.22%  "
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("E13").Value = "  -5.37%  "
$ws.Range("D14").Value = "2.498.87"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.91"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("E16").Value = "  -4.32%  "
$ws.Range("D17").Value = "2.170.79"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "41.014.61"
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.37"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.06"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "225.28"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -7.62%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.83"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.53"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.79"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.94"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "30.82"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0769"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.15"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -8.80%  "
$ws.Range("E35").Value = "  -3.19%  "
$ws.Range("E36").Value = "  -9.03%  "
$ws.Range("E37").Value = "  -3.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0285"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.25"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -5.91%  "
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("E41").Value = "  -4.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "60.11"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -8.06%  "
$ws.Range("E43").Value = "  -4.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.32"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0971"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.15"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.96%  "
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("E49").Value = "  -7.66%  "
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("D51").Value = "2.376.76"
$ws.Range("E51").Value = "  -1.96%  "
